# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data currently spans A1:AC53).
$lastRow = $ws.UsedRange.Rows.Count

# New header cells in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style already used for the other header cells (bold, bordered,
# centered/top-aligned) by copying it from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (67 wins, 94 losses, 0 ties) for every player row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 94
    $ws.Cells.Item($r, 32).Value = 0
}
